# "Roo is up and running" — add two new dictionary entries (chico/-a, llegar)
# and reformat the "part_of_speech" (column D) cells to show their value
# wrapped in parentheses, e.g. "(Verb)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) New custom number format "(@)" applied (with wrap text) to the existing
#    non-blank part_of_speech cells in column D (rows 2,6,7,8,9). This is
#    what produces the new numFmt/cellXf (style index 10) in styles.xml.
# ---------------------------------------------------------------------------
$dCells = $ws.Range("D2:D9")
$dCells.NumberFormat = "\(@\)"
$dCells.WrapText = $true

# The three blank/empty styled cells in column D (rows 3,4,5) are fully
# cleared instead of just reformatted, so they drop out of the sheet XML.
$ws.Range("D3").Clear()
$ws.Range("D4").Clear()
$ws.Range("D5").Clear()

# ---------------------------------------------------------------------------
# 2) New row 10: chico, -a
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "chico, -a"
$ws.Range("D10").Value = "Adjective, Noun"
$ws.Range("E10").Value = 'small; (masc.) "boy," "child;" (fem.) "girl"'
$ws.Range("F10").Value = 'Perhaps from Latin ciccus "nothing," earlier "something worthless," but originally "the thin membrane surrounding the grains of a pomegranate." Presumably borrowed from an unattested Ancient Greek word *κίκκος (kíkkos) "shell of a pomegranate," hypothesized by Beekes (2008) on the basis of the Latin word and possible Greek derivatives κίκκαβος (kíkkabos) "small coin in the Underworld," κικκάβι(ο)ν (kikkábi(o)n) "nothing," and κικαῖος (kikaios) of obscure meaning. '
$ws.Range("G10").Value = "Ultimately of unknown origin."
$ws.Range("H10").Value = 'The sense of "small" was first and then was extended to children. The change from c- to ch- in Spanish and Asturian is unexpected and may be due to Basque influence, via  txiki "small," "few," from earlier tiki. However, the change is not unique (e.g. Latin cicer "chick-pea" > Spanish chícharo "pea").'

$ws.Range("D10").NumberFormat = "\(@\)"
$ws.Range("D10").WrapText = $true

$ws.Range("A10:H10").Font.Name = "Georgia"
$ws.Range("A10:H10").Font.Size = 10
$ws.Range("A10").WrapText = $true
$ws.Range("E10:H10").WrapText = $true

$ws.Rows.Item(10).RowHeight = 79

# ---------------------------------------------------------------------------
# 3) New row 11: llegar
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "llegar"
$ws.Range("D11").Value = "Verb"
$ws.Range("E11").Value = "to arrive"
$ws.Range("F11").Value = "From Latin plicare ""to fold."" According to Roberts (2014), an extension of the sense of being 'folded' into something as an arrival at a destination."
$ws.Range("G11").Value = "From Proto-Italic *plek- 'id.' From Proto-Indo-European *pleḱ- 'id.' "

$ws.Range("D11").NumberFormat = "\(@\)"
$ws.Range("D11").WrapText = $true

$ws.Range("A11:G11").Font.Name = "Georgia"
$ws.Range("A11:G11").Font.Size = 10
$ws.Range("A11").WrapText = $true
$ws.Range("E11:G11").WrapText = $true

$ws.Rows.Item(11).RowHeight = 51

# ---------------------------------------------------------------------------
# 4) View state: scrolled down so row 8 is the first visible row below the
#    frozen header, with E9 as the active cell.
# ---------------------------------------------------------------------------
$ws.Range("E9").Select()
$excel.ActiveWindow.ScrollRow = 8
